# Adiciona as colunas "fonte" e "observacao" na planilha de upload.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cabecalhos das novas colunas
$ws.Range("H1").Value = "fonte"
$ws.Range("I1").Value = "observacao"

# Linha 2
$ws.Range("H2").Value = "isso eh uma fonte"

# Linha 3 fica sem valores nas novas colunas (apenas formatadas)

# Linha 4
$ws.Range("H4").Value = "segunda fonte"
$ws.Range("I4").Value = "alguma obs"

# Aplica o novo estilo (fonte normal aplicada explicitamente) as novas celulas
$ws.Range("H1:I4").Font.Name = "Arial"
$ws.Range("H1:I4").Font.Size = 10

# Seleciona a area recem adicionada, como ficou no arquivo final
$ws.Range("H1:I4").Select()
